$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.423.19'
$ws.Range("E2").Value = '  +2.09%  '
$ws.Range("D3").Value = '1.845.29'
$ws.Range("E3").Value = '  +1.82%  '
$ws.Range("E4").Value = '  +1.36%  '
$ws.Range("E5").Value = '  +2.42%  '
$ws.Range("D6").Value = "'1.014"
$ws.Range("E6").Value = '  +1.23%  '
$ws.Range("E7").Value = '  +1.98%  '
$ws.Range("D8").Value = "'0.3705"
$ws.Range("E8").Value = '  +0.60%  '
$ws.Range("D9").Value = "'0.07459"
$ws.Range("E9").Value = '  +1.37%  '
$ws.Range("D10").Value = "'0.8879"
$ws.Range("E10").Value = '  +2.29%  '
$ws.Range("D11").Value = "'20.53"
$ws.Range("E11").Value = '  +0.75%  '
$ws.Range("D12").Value = '1.855.06'
$ws.Range("E12").Value = '  +5.53%  '
$ws.Range("D13").Value = "'0.07415"
$ws.Range("E13").Value = '  +4.99%  '
$ws.Range("D14").Value = "'5.492"
$ws.Range("E14").Value = '  +2.84%  '
$ws.Range("D15").Value = "'93.34"
$ws.Range("E15").Value = '  +1.72%  '
$ws.Range("D16").Value = "'6.589"
$ws.Range("E16").Value = '  +1.62%  '
$ws.Range("E17").Value = '  +1.26%  '
$ws.Range("D18").Value = "'0.000008863"
$ws.Range("E18").Value = '  +2.04%  '
$ws.Range("E19").Value = '  +1.30%  '
$ws.Range("D20").Value = "'14.87"
$ws.Range("E20").Value = '  +1.00%  '
$ws.Range("D21").Value = '27.447.23'
$ws.Range("E21").Value = '  +2.11%  '
$ws.Range("D22").Value = "'5.340"
$ws.Range("E22").Value = '  +0.28%  '
$ws.Range("D23").Value = "'10.72"
$ws.Range("E23").Value = '  +1.76%  '
$ws.Range("D24").Value = '2.074.70'
$ws.Range("E24").Value = '  +3.54%  '
$ws.Range("D25").Value = "'1.913"
$ws.Range("E25").Value = '  +0.77%  '
$ws.Range("D26").Value = "'152.41"
$ws.Range("E26").Value = '  +0.94%  '
$ws.Range("D27").Value = "'18.63"
$ws.Range("E27").Value = '  +1.34%  '
$ws.Range("D28").Value = "'2.178"
$ws.Range("E28").Value = '  +0.89%  '
$ws.Range("D29").Value = "'5.293"
$ws.Range("E29").Value = '  -0.24%  '
$ws.Range("D30").Value = "'118.15"
$ws.Range("E30").Value = '  +2.23%  '
$ws.Range("D31").Value = "'0.08983"
$ws.Range("E31").Value = '  +0.67%  '
$ws.Range("D32").Value = "'0.7615"
$ws.Range("E32").Value = '  -0.39%  '
$ws.Range("D33").Value = "'1.179"
$ws.Range("E33").Value = '  +2.09%  '
$ws.Range("D34").Value = "'4.568"
$ws.Range("E34").Value = '  +1.61%  '
$ws.Range("D35").Value = "'2.953"
$ws.Range("E35").Value = '  +1.69%  '
$ws.Range("E36").Value = '  +1.30%  '
$ws.Range("D37").Value = "'1.110"
$ws.Range("E37").Value = '  +2.03%  '
$ws.Range("D38").Value = "'0.05367"
$ws.Range("E38").Value = '  +1.62%  '
$ws.Range("E39").Value = '  +0.59%  '
$ws.Range("D40").Value = "'3.012"
$ws.Range("E40").Value = '  +2.46%  '
$ws.Range("D41").Value = "'7.322"
$ws.Range("E41").Value = '  +1.01%  '
$ws.Range("D42").Value = "'2.392"
$ws.Range("E42").Value = '  +2.01%  '
$ws.Range("D43").Value = "'0.5353"
$ws.Range("E43").Value = '  +1.10%  '
$ws.Range("D44").Value = "'0.1668"
$ws.Range("E44").Value = '  +0.35%  '
$ws.Range("D45").Value = "'8.543"
$ws.Range("E45").Value = '  +1.76%  '
$ws.Range("D46").Value = "'0.4967"
$ws.Range("E46").Value = '  +1.01%  '
$ws.Range("D47").Value = "'10.57"
$ws.Range("E47").Value = '  +1.73%  '
$ws.Range("D48").Value = "'1.015"
$ws.Range("E48").Value = '  +1.40%  '
$ws.Range("D49").Value = "'105.06"
$ws.Range("E49").Value = '  +1.66%  '
$ws.Range("D50").Value = "'1.685"
$ws.Range("E50").Value = '  +1.19%  '
$ws.Range("D51").Value = "'0.06334"
$ws.Range("E51").Value = '  +0.89%  '
